# Auto-generated edit script: update crypto price/volume data per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells whose new value is plain text (not interpretable as a number) ---
# Safe to assign directly.
$ws.Range('D2').Value = '69.823.45'
$ws.Range('E2').Value = '  -1.20%  '
$ws.Range('D3').Value = '3.503.98'
$ws.Range('E3').Value = '  -1.67%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('E5').Value = '  +3.34%  '
$ws.Range('E6').Value = '  +1.29%  '
$ws.Range('E7').Value = '  +0.70%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('E9').Value = '  -1.24%  '
$ws.Range('E10').Value = '  +2.67%  '
$ws.Range('E11').Value = '  -1.26%  '
$ws.Range('E12').Value = '  -1.12%  '
$ws.Range('E13').Value = '  +2.27%  '
$ws.Range('D14').Value = '4.063.62'
$ws.Range('E14').Value = '  -1.58%  '
$ws.Range('E15').Value = '  +8.74%  '
$ws.Range('D16').Value = '69.891.89'
$ws.Range('E16').Value = '  -1.10%  '
$ws.Range('E17').Value = '  +0.00%  '
$ws.Range('E18').Value = '  -0.29%  '
$ws.Range('D19').Value = '3.501.57'
$ws.Range('E19').Value = '  -2.06%  '
$ws.Range('E20').Value = '  -0.18%  '
$ws.Range('E21').Value = '  -0.48%  '
$ws.Range('E22').Value = '  +0.13%  '
$ws.Range('E23').Value = '  +12.42%  '
$ws.Range('E24').Value = '  +0.43%  '
$ws.Range('E25').Value = '  +3.11%  '
$ws.Range('E26').Value = '  +4.08%  '
$ws.Range('E27').Value = '  -0.54%  '
$ws.Range('E28').Value = '  +5.85%  '
$ws.Range('E30').Value = '  +0.73%  '
$ws.Range('B31').Value = 'dogwifhat'
$ws.Range('C31').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('E31').Value = '  +5.59%  '
$ws.Range('B32').Value = 'Cosmos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('E32').Value = '  +3.99%  '
$ws.Range('E33').Value = '  +0.16%  '
$ws.Range('E34').Value = '  +1.47%  '
$ws.Range('D35').Value = '3.738.88'
$ws.Range('E35').Value = '  +1.82%  '
$ws.Range('B36').Value = 'Bittensor'
$ws.Range('C36').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('E36').Value = '  -1.03%  '
$ws.Range('B37').Value = 'Fetch.AI'
$ws.Range('C37').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('E37').Value = '  -4.35%  '
$ws.Range('E39').Value = '  +1.26%  '
$ws.Range('B40').Value = 'InjectiveProtocol'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('E40').Value = '  -3.62%  '
$ws.Range('B41').Value = 'TheGraph'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('E41').Value = '  -3.75%  '
$ws.Range('E42').Value = '  +0.14%  '
$ws.Range('E43').Value = '  +0.17%  '
$ws.Range('E44').Value = '  +1.04%  '
$ws.Range('B45').Value = 'ThetaToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('E45').Value = '  -2.64%  '
$ws.Range('B46').Value = 'Stellar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('E46').Value = '  +2.36%  '
$ws.Range('E47').Value = '  -4.74%  '
$ws.Range('B48').Value = 'FirstDigitalUSD'
$ws.Range('C48').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('E48').Value = '  +0.45%  '
$ws.Range('B49').Value = 'THORChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('E49').Value = '  -4.71%  '
$ws.Range('E50').Value = '  -1.32%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('E51').Value = '  +11.07%  '

# --- Cells whose new value LOOKS numeric (e.g. "605.91", "1.00") but must
#     remain stored as TEXT, matching the column's existing inline-string type.
#     Force text via NumberFormat "@" before assignment, then clear the
#     number-format override (ClearFormats) so no stray per-cell style is
#     left behind, while the value itself stays text.
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '605.91'
$ws.Range('D5').ClearFormats()
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '191.74'
$ws.Range('D6').ClearFormats()
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.662'
$ws.Range('D10').ClearFormats()
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '53.47'
$ws.Range('D11').ClearFormats()
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '617.79'
$ws.Range('D15').ClearFormats()
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.991'
$ws.Range('D21').ClearFormats()
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '17.99'
$ws.Range('D22').ClearFormats()
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '105.53'
$ws.Range('D23').ClearFormats()
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.01'
$ws.Range('D25').ClearFormats()
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.91'
$ws.Range('D28').ClearFormats()
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '34.17'
$ws.Range('D29').ClearFormats()
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.19'
$ws.Range('D31').ClearFormats()
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '12.67'
$ws.Range('D32').ClearFormats()
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '64.11'
$ws.Range('D34').ClearFormats()
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '524.79'
$ws.Range('D36').ClearFormats()
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.09'
$ws.Range('D37').ClearFormats()
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.00'
$ws.Range('D38').ClearFormats()
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '36.80'
$ws.Range('D40').ClearFormats()
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.390'
$ws.Range('D41').ClearFormats()
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0463'
$ws.Range('D44').ClearFormats()
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.86'
$ws.Range('D45').ClearFormats()
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.141'
$ws.Range('D46').ClearFormats()
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.32'
$ws.Range('D47').ClearFormats()
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.00'
$ws.Range('D48').ClearFormats()
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.74'
$ws.Range('D49').ClearFormats()
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '132.60'
$ws.Range('D50').ClearFormats()
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.30'
$ws.Range('D51').ClearFormats()

Write-Host "Applied 103 cell updates (76 text + 27 numeric-like)"
